# RPA datasets push 2023-10-18
#
# The "그린리소스" (GreenResource) IPO listing's demand-forecast date moved
# up from 2023.10.25~10.31 to 2023.11.03~11.09, so its row needs to move
# from its old sorted position (row 12, right after 한국스팩13호) up to
# right after 와이바이오로직스 (row 2) to keep the table sorted by date.
#
# All other figures for this listing (price band, offering amount,
# underwriter) stay the same - only its position and date change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the row at its new (earlier) position.
$ws.Rows(3).Insert()

# Populate the newly-inserted row with 그린리소스's data, using the
# corrected/updated demand-forecast date.
$ws.Range("A3").Value = "그린리소스"
$ws.Range("B3").Value = "2023.11.03~11.09"
$ws.Range("C3").Value = "11,000~14,000"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 18040
$ws.Range("F3").Value = "NH투자증권"

# Remove the old 그린리소스 row, which the insert above shifted down to
# row 13.
$ws.Rows(13).Delete()
